# Resort the sheet tabs: reverse their order.
# Before: 2020-Q4, 2021-Q2, 2021-Q3, 2022-Q1, 2022-Q2, 总计
# After:  总计, 2022-Q2, 2022-Q1, 2021-Q3, 2021-Q2, 2020-Q4
$wb = $excel.ActiveWorkbook

$count = $wb.Worksheets.Count

# Capture the current left-to-right tab order first (names are stable
# identifiers even as positions change below).
$names = @()
for ($i = 1; $i -le $count; $i++) {
    $names += $wb.Worksheets.Item($i).Name
}

# Walking the original order and moving each sheet in turn to the very
# front reverses the overall tab order.
foreach ($n in $names) {
    $sheet = $wb.Worksheets.Item($n)
    $sheet.Move($wb.Worksheets.Item(1))
}
